$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for the new columns.
#    Old layout:  A B C D E F G H I
#      G = Tempo di esecuzione (s)
#      H = Precisione media (%)
#      I = Precisione (%)
#    New layout:  A B C D E F G H I J K L M
#      G = Tempo estrazione (s)          [new, empty]
#      H = Tempo calcolo distanze (s)    [new, empty]
#      I = Tempo di esecuzione (s)       [old G]
#      J = Numero cluster corretti       [new, empty]
#      K = Numero elementi corretti      [new, empty]
#      L = Precisione media (%)          [old H]
#      M = Precisione (%)                [old I]
# ---------------------------------------------------------------------------
$ws.Range("G:H").Insert()
$ws.Range("J:K").Insert()

# ---------------------------------------------------------------------------
# 2. Header row (row 2) text updates.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Numero scansioni"
$ws.Range("C2").Value = "Numero di elementi stimato"
$ws.Range("G2").Value = "Tempo estrazione (s)"
$ws.Range("H2").Value = "Tempo calcolo distanze (s)"
$ws.Range("I2").Value = "Tempo di esecuzione (s)"
$ws.Range("J2").Value = "Numero cluster corretti"
$ws.Range("K2").Value = "Numero elementi corretti"
$ws.Range("L2").Value = "Precisione media (%)"
$ws.Range("M2").Value = "Precisione (%)"

# ---------------------------------------------------------------------------
# 3. Row labels (column A) renamed from "immagini" to "scansioni" for the
#    700-count rows; the 200/100 rows keep their existing text.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "700 scansioni (LCS e L1)"
$ws.Range("A4").Value = "700 scansioni (L1)"
$ws.Range("A5").Value = "700 scansioni (LCS)"
$ws.Range("A6").Value = "200 immagini (LCS)"
$ws.Range("A7").Value = "200 immagini (LCS e L1)"
$ws.Range("A8").Value = "100 immagini (LCS)"
$ws.Range("A9").Value = "100 immagini (LCS e L1)"

# ---------------------------------------------------------------------------
# 4. Updated data point: M3 (Precisione (%) for the 700 LCS e L1 row).
# ---------------------------------------------------------------------------
$ws.Range("M3").Value = 76.3

# ---------------------------------------------------------------------------
# 5. Column widths (best match on the engine's 1/6-character-width grid).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.666666666666668
$ws.Columns.Item(2).ColumnWidth = 17.833333333333332
$ws.Columns.Item(3).ColumnWidth = 25.0
$ws.Columns.Item(4).ColumnWidth = 23.0
$ws.Columns.Item(5).ColumnWidth = 14.5
$ws.Columns.Item(6).ColumnWidth = 22.833333333333332
$ws.Columns.Item(7).ColumnWidth = 19.333333333333332
$ws.Columns.Item(8).ColumnWidth = 23.5
$ws.Columns.Item(9).ColumnWidth = 22.5
$ws.Columns.Item(10).ColumnWidth = 21.333333333333332
$ws.Columns.Item(11).ColumnWidth = 23.0
$ws.Columns.Item(12).ColumnWidth = 18.333333333333332
$ws.Columns.Item(13).ColumnWidth = 12.5

# ---------------------------------------------------------------------------
# 6. Selection matches the author's last-saved cursor position.
# ---------------------------------------------------------------------------
[void]$ws.Range("L12").Select()

Write-Host "done"
